$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old totals row (row 10) so the table grows
# from 9 data rows to 10 data rows, pushing the totals row from 10 -> 11.
$ws.Rows.Item(10).Insert()

# Updated band-edge data (DiMascolo25 band edges), rows 2-10.
$data = @(
    @(30, 54, 4, 40),
    @(66, 117, 4, 120),
    @(120, 182, 4, 120),
    @(183, 252, 8, 250),
    @(252, 325, 8, 250),
    @(325, 375, 18, 375),
    @(384, 422, 18, 375),
    @(595, 713, 10, 640),
    @(786, 905, 6, 870)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Totals row now lives at row 11.
$ws.Range("C11").Formula = "=SUM(C2,C3,C5,C7,C9,C10)"

# Selection moved to C9 in the saved file.
$ws.Range("C9").Select()
